$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7154039144515991
$ws.Range("B1").Value = 1.489919304847717
$ws.Range("C1").Value = 4.086187839508057
$ws.Range("D1").Value = 2.648882389068604
$ws.Range("E1").Value = 0.8133031725883484
